$d = $word.ActiveDocument

$replacements = @(
    @{old="200×7=1400"; new="757×6=4542"},
    @{old="515×5=2575"; new="692×6=4152"},
    @{old="487×9=4383"; new="895×5=4475"},
    @{old="879×4=3516"; new="621×4=2484"},
    @{old="276×2=552"; new="630×6=3780"},
    @{old="947×8=7576"; new="433×9=3897"},
    @{old="624×2=1248"; new="611×5=3055"},
    @{old="378×7=2646"; new="227×5=1135"},
    @{old="811×9=7299"; new="198×4=792"},
    @{old="795×3=2385"; new="639×6=3834"},
    @{old="708×4=2832"; new="800×2=1600"},
    @{old="249×4=996"; new="841×2=1682"},
    @{old="210×5=1050"; new="624×6=3744"},
    @{old="177×7=1239"; new="184×3=552"},
    @{old="897×9=8073"; new="981×9=8829"},
    @{old="135×2=270"; new="813×5=4065"},
    @{old="887×6=5322"; new="161×6=966"},
    @{old="369×5=1845"; new="793×2=1586"},
    @{old="192×7=1344"; new="575×5=2875"},
    @{old="186×8=1488"; new="593×3=1779"},
    @{old="635×8=5080"; new="401×3=1203"},
    @{old="484×8=3872"; new="741×5=3705"},
    @{old="665×5=3325"; new="500×5=2500"},
    @{old="367×5=1835"; new="756×7=5292"},
    @{old="372×7=2604"; new="932×9=8388"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
